$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with the new data-driven test values
$ws.Range("A2").Value = "uname1"
$ws.Range("A3").Value = "uname2"
$ws.Range("A4").Value = "uname3"
$ws.Range("A5").Value = "uname4"

$ws.Range("B2").Value = "password1"
$ws.Range("B3").Value = "password2"
$ws.Range("B4").Value = "password3"
$ws.Range("B5").Value = "password4"

$ws.Range("D2").Value = "lastname1"
$ws.Range("D3").Value = "lastname2"
$ws.Range("D4").Value = "lastname3"
$ws.Range("D5").Value = "lastname4"

$ws.Range("E2").Value = "firstname1"
$ws.Range("E3").Value = "firstname2"
$ws.Range("E4").Value = "firstname3"
$ws.Range("E5").Value = "firstname4"

# Move active selection to G7
$ws.Range("G7").Select() | Out-Null
